# Update countries & provincias Spain
# Applies the 17-Aug-2020 data refresh (07:33 -> 08:50) to the "Pais" sheet:
#  - updates case totals for several countries
#  - re-sorts a few adjacent rows whose updated totals changed their rank
#    (Ucrania overtakes Israel, Georgia overtakes Republica de Chipre,
#     Islas Malvinas overtakes Montserrat)
#  - updates the "Datos actualizados..." timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: India (no reordering, just refreshed totals) ---
$ws.Range("B6").Value = 2648353
$ws.Range("C6").Value = 1037
$ws.Range("D6").Value = 1920217
$ws.Range("E6").Value = 677077
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 14
$ws.Range("H6").Value = 51059

# --- Rows 33-34: Ucrania overtakes Israel ---
# Row 33 becomes Ucrania with its freshly updated totals.
$ws.Range("A33").Value = "Ucrania"
$ws.Range("B33").Value = 92820
$ws.Range("C33").Value = 1464
$ws.Range("D33").Value = 48164
$ws.Range("E33").Value = 42567
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 19
$ws.Range("H33").Value = 2089

# Row 34 becomes Israel, keeping its previous (unchanged) totals.
$ws.Range("A34").Value = "Israel"
$ws.Range("B34").Value = 92680
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 68510
$ws.Range("E34").Value = 23485
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 685

# --- Row 73: El Salvador (refreshed totals) ---
$ws.Range("D73").Value = 10814
$ws.Range("E73").Value = 11480
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 6
$ws.Range("H73").Value = 618

# --- Row 108: Hungria (refreshed totals) ---
$ws.Range("B108").Value = 4946
$ws.Range("C108").Value = 30
$ws.Range("D108").Value = 3630
$ws.Range("E108").Value = 708

# --- Rows 145-146: Georgia overtakes Republica de Chipre ---
# Row 145 becomes Georgia with its freshly updated totals.
$ws.Range("A145").Value = "Georgia"
$ws.Range("B145").Value = 1341
$ws.Range("C145").Value = 5
$ws.Range("D145").Value = 1092
$ws.Range("E145").Value = 232
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 17

# Row 146 becomes Republica de Chipre, keeping its previous (unchanged) totals.
$ws.Range("A146").Value = "Republica de Chipre"
$ws.Range("B146").Value = 1339
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 870
$ws.Range("E146").Value = 449
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 20

# --- Row 147: Letonia (refreshed totals) ---
$ws.Range("B147").Value = 1323
$ws.Range("C147").Value = 1
$ws.Range("E147").Value = 213

# --- Rows 213-214: Islas Malvinas overtakes Montserrat ---
# Row 213 becomes Islas Malvinas.
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 13
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

# Row 214 becomes Montserrat.
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1

# --- Timestamp string update ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Agosto de 2020 a las 08:50"
